# Reorders header columns K:AM (the "jobgroup" column, which was sitting at the
# end of the salary block in AM, is moved so it becomes the first salary-block
# column at K) for both the label row (row 1) and the description row (row 2),
# then updates the sheet's view/selection.
#
# This mirrors an Excel "cut column AM, insert before column K" operation:
# the old AM value moves to K, and everything that used to occupy K..AL shifts
# one column to the right to become L..AM. Columns before K and after AM are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startCol = 11  # column K
$endCol   = 39  # column AM

# Capture the current values across the block for the two populated rows.
$headerRow = @()
$descRow   = @()
for ($c = $startCol; $c -le $endCol; $c++) {
    $headerRow += ,$ws.Cells.Item(1, $c).Value2
    $descRow   += ,$ws.Cells.Item(2, $c).Value2
}

# Rotate: last column's value becomes the first, all others shift right by one.
$count = $headerRow.Count
$newHeaderRow = @($headerRow[$count - 1]) + $headerRow[0..($count - 2)]
$newDescRow   = @($descRow[$count - 1]) + $descRow[0..($count - 2)]

$i = 0
for ($c = $startCol; $c -le $endCol; $c++) {
    $ws.Cells.Item(1, $c).Value2 = $newHeaderRow[$i]
    $ws.Cells.Item(2, $c).Value2 = $newDescRow[$i]
    $i++
}

# Reset the view: scroll back to the top-left (no frozen/offset topLeftCell)
# and move the selection to K2.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K2").Select() | Out-Null
